$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 4 and row 5 for columns D, J, K, M, P
$d4 = $ws.Range("D4").Value2
$j4 = $ws.Range("J4").Value2
$k4 = $ws.Range("K4").Value2
$m4 = $ws.Range("M4").Value2
$p4 = $ws.Range("P4").Value2

$d5 = $ws.Range("D5").Value2
$j5 = $ws.Range("J5").Value2
$k5 = $ws.Range("K5").Value2
$m5 = $ws.Range("M5").Value2
$p5 = $ws.Range("P5").Value2

$ws.Range("D4").Value2 = $d5
$ws.Range("J4").Value2 = $j5
$ws.Range("K4").Value2 = $k5
$ws.Range("M4").Value2 = $m5
$ws.Range("P4").Value2 = $p5

$ws.Range("D5").Value2 = $d4
$ws.Range("J5").Value2 = $j4
$ws.Range("K5").Value2 = $k4
$ws.Range("M5").Value2 = $m4
$ws.Range("P5").Value2 = $p4
